# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (want-to-go count) figures in the 展览 sheet and
# the mirrored rows in the 全部类型 aggregate sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 830
$ws1.Range("F13").Value = 13581
$ws1.Range("F17").Value = 5578
$ws1.Range("F18").Value = 5590
$ws1.Range("F19").Value = 65

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 830
$ws4.Range("F35").Value = 13581
$ws4.Range("F40").Value = 5578
$ws4.Range("F41").Value = 5590
$ws4.Range("F42").Value = 65
